# AD-EYE/TA/SimulinkConfig.xlsx — bugfixes around folder/experiment names.
# Append two new BlockName/Value rows to Sheet1 and update the sheet's
# view state (scroll position + active selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows appended at the bottom of the table.
$ws.Range("A29").Value = "percent_reflecting_sfc"
$ws.Range("B29").Value = 0.9

$ws.Range("A30").Value = "R"
$ws.Range("B30").Value = 0

# Bring the sheet into view and scroll so row 7 is at the top, then
# leave the selection on the first empty cell below the new data (B31),
# matching where the user's cursor ended up after the edit.
$ws.Select()
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
